$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 69-72 had their player name (column A) and weekly score columns
# (E-H) cyclically shifted up by one row: row69 <- row70, row70 <- row71,
# row71 <- row72, row72 <- row69 (wrap-around).

# Column A (player names)
$ws.Cells.Item(69, 1).Value = "DGJ-DAVI"
$ws.Cells.Item(70, 1).Value = "luck"
$ws.Cells.Item(71, 1).Value = "andrebts"
$ws.Cells.Item(72, 1).Value = "Asten Acady"

# Row 69 (E,F,G,H)
$ws.Cells.Item(69, 5).Value = 16
$ws.Cells.Item(69, 6).Value = 16
$ws.Cells.Item(69, 7).Value = 12
$ws.Cells.Item(69, 8).Value = 15

# Row 70 (E,F,G,H)
$ws.Cells.Item(70, 5).Value = 16
$ws.Cells.Item(70, 6).Value = 16
$ws.Cells.Item(70, 7).Value = 12
$ws.Cells.Item(70, 8).Value = 16

# Row 71 (E,F,G,H)
$ws.Cells.Item(71, 5).Value = 15
$ws.Cells.Item(71, 6).Value = 16
$ws.Cells.Item(71, 7).Value = 16
$ws.Cells.Item(71, 8).Value = 16

# Row 72 (E,F,G,H)
$ws.Cells.Item(72, 5).Value = 16
$ws.Cells.Item(72, 6).Value = 16
$ws.Cells.Item(72, 7).Value = 16
$ws.Cells.Item(72, 8).Value = 16
